$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: G3 and H3 -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: D4 and E4 -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: D5 and E5 -> 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: D6 and E6 -> 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7: H7 -> 1
$ws.Range("H7").Value = 1

# Row 8: H8 -> 1
$ws.Range("H8").Value = 1

# Row 9: H9 -> 1
$ws.Range("H9").Value = 1

# Row 10: H10 -> 1
$ws.Range("H10").Value = 1

# Row 11: H11 -> 1
$ws.Range("H11").Value = 1

# Row 12: D12 and E12 -> 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: H13 -> 1
$ws.Range("H13").Value = 1

# Row 14: H14 -> 1
$ws.Range("H14").Value = 1

# Row 15: H15 -> 1
$ws.Range("H15").Value = 1

# Row 16: H16 -> 1
$ws.Range("H16").Value = 1

# Row 17: H17 -> 1
$ws.Range("H17").Value = 1

# Row 18: H18 -> 1
$ws.Range("H18").Value = 1
